$d = $word.ActiveDocument

# Update the date/title line at the top of the document
$d.Content.Find.Execute("2025-08-10 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-11 Monday", 2)

# Update the division problems in the table, using row/cell indices to
# avoid ambiguity from duplicate text (e.g. "13÷7=" appears twice).
$table = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="52÷2="},
    @{Row=1;  Col=2; New="54÷5="},
    @{Row=1;  Col=3; New="25÷3="},
    @{Row=1;  Col=4; New="22÷4="},
    @{Row=1;  Col=5; New="40÷9="},

    @{Row=5;  Col=1; New="26÷8="},
    @{Row=5;  Col=2; New="41÷5="},
    @{Row=5;  Col=3; New="40÷5="},
    @{Row=5;  Col=4; New="71÷9="},
    @{Row=5;  Col=5; New="27÷4="},

    @{Row=9;  Col=1; New="36÷7="},
    @{Row=9;  Col=2; New="16÷2="},
    @{Row=9;  Col=3; New="98÷3="},
    @{Row=9;  Col=4; New="78÷7="},
    @{Row=9;  Col=5; New="29÷5="},

    @{Row=13; Col=1; New="74÷9="},
    @{Row=13; Col=2; New="16÷8="},
    @{Row=13; Col=3; New="26÷5="},
    @{Row=13; Col=4; New="11÷6="},
    @{Row=13; Col=5; New="67÷6="},

    @{Row=17; Col=1; New="21÷4="},
    @{Row=17; Col=2; New="41÷4="},
    @{Row=17; Col=3; New="11÷7="},
    @{Row=17; Col=4; New="26÷9="},
    @{Row=17; Col=5; New="14÷3="}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $range = $cell.Range
    # Trim the trailing cell-mark/paragraph-mark characters so we only
    # replace the visible text content.
    $range.End = $range.End - 1
    $range.Text = $u.New
}
